$wb = $excel.ActiveWorkbook

# --- Update the daily conversion text on sheet "Hoja1" ---
$ws1 = $wb.Worksheets.Item("Hoja1")
$cell = $ws1.Range("A1")
$text = $cell.Value2
$text = $text.Replace("1000 Bs = 3.14 = 11772.05 pesos", "1000 Bs = 3.15 = 11810.32 pesos")
$text = $text.Replace("11772.05 pesos = 3.12 = 970.42 Bs", "11810.32 pesos = 3.14 = 973.72 Bs")
$cell.Value2 = $text

# --- Update the rate values on sheet "tasas" ---
$ws2 = $wb.Worksheets.Item("tasas")
$ws2.Range("N10").Value = 317
$ws2.Range("O10").Value = 3743.87
$ws2.Range("N12").Value = 3760
$ws2.Range("O12").Value = 310
